$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: a new 2022-Q3 row is published, so every existing
#    quarter shifts down by one row. Insert a fresh row 2, fill it with the
#    2022-Q3 totals, and renumber the running index in column A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 2.97

# Column A keeps the bordered/bold style used by the rest of the index
# column; copy it down from the row just below instead of re-deriving it.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

for ($r = 3; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) Insert the brand-new "2022-Q3" worksheet right after "总计" (pushing the
#    older quarter tabs one position later, matching the workbook.xml diff).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Match the header look (style) used on every other quarter sheet by copying
# the formatting from the "总计" header cell, then overwrite the text.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows: column A is a numeric running index sharing the "总计"-style
# index-column formatting; B-G are free-form text (fund code/name/percentages
# kept as strings, not converted to numbers); H is a plain number.
$total.Cells.Item(2, 1).Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)

$rows = @(
    @(0, "900090", "中信卓越成长两年持有期混合B", "51.50", "93.14", "3.78", "1.9467", 7),
    @(1, "900010", "中信卓越成长两年持有期混合A", "14.24", "93.14", "3.78", "0.5383", 7),
    @(2, "003396", "东方红优享红利混合",           "14.04", "60.52", "2.23", "0.3131", 9),
    @(3, "900100", "中信卓越成长两年持有期混合C", "4.61",  "93.14", "3.78", "0.1743", 7)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q3.Cells.Item($r, 1).Value = $data[0]

    for ($c = 1; $c -le 6; $c++) {
        $cell = $q3.Cells.Item($r, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c]
        $cell.ClearFormats()
    }

    $q3.Cells.Item($r, 8).Value = $data[7]
}

Write-Output "done"
